$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Nombre"
$ws.Range("B1").Value = "Matrícula"
# Prefix with an apostrophe so Excel stores this as literal text
# instead of auto-converting it to a date serial value.
$ws.Range("C1").Value = "'2025-05-24"

$ws.Range("A2").Value = "Ignacio Nava Castillo"
$ws.Range("B2").Value = "PG2001"
